$wb = $excel.ActiveWorkbook

# 1) Remove the "Location" sheet entirely
$wb.Worksheets("Location").Delete()

# 2) Populate new data rows (2-12) on "Incomplete Units"
$ws6 = $wb.Worksheets("Incomplete Units")
$ws6.Range("A2").Value = 0
$ws6.Range("B2").Value = 0
$ws6.Range("C2").Value = 0
$ws6.Range("D2").Value = 0
$ws6.Range("E2").Value = 0
$ws6.Range("F2").Value = 0
$ws6.Range("G2").Value = 0
$ws6.Range("H2").Value = 0
$ws6.Range("I2").Value = 0
$ws6.Range("J2").Value = 2
$ws6.Range("K2").Value = 4
$ws6.Range("L2").Value = 1028
$ws6.Range("N2").Value = 0
$ws6.Range("O2").Value = 0
$ws6.Range("P2").Value = 0
$ws6.Range("Q2").Value = 9
$ws6.Range("R2").Value = 16
$ws6.Range("A3").Value = 0
$ws6.Range("B3").Value = 0
$ws6.Range("C3").Value = 0
$ws6.Range("D3").Value = 0
$ws6.Range("E3").Value = 0
$ws6.Range("F3").Value = 0
$ws6.Range("G3").Value = 0
$ws6.Range("H3").Value = 0
$ws6.Range("I3").Value = 0
$ws6.Range("J3").Value = 1
$ws6.Range("K3").Value = 0
$ws6.Range("N3").Value = 0
$ws6.Range("O3").Value = 0
$ws6.Range("P3").Value = 0
$ws6.Range("Q3").Value = 0
$ws6.Range("R3").Value = 7
$ws6.Range("A4").Value = 0
$ws6.Range("B4").Value = 0
$ws6.Range("C4").Value = 0
$ws6.Range("D4").Value = 0
$ws6.Range("E4").Value = 0
$ws6.Range("F4").Value = 0
$ws6.Range("G4").Value = 0
$ws6.Range("H4").Value = 1
$ws6.Range("I4").Value = 0
$ws6.Range("J4").Value = 0
$ws6.Range("K4").Value = 263
$ws6.Range("N4").Value = 0
$ws6.Range("O4").Value = 0
$ws6.Range("P4").Value = 0
$ws6.Range("Q4").Value = 0
$ws6.Range("R4").Value = 0
$ws6.Range("A5").Value = 0
$ws6.Range("B5").Value = 0
$ws6.Range("C5").Value = 0
$ws6.Range("D5").Value = 0
$ws6.Range("E5").Value = 0
$ws6.Range("F5").Value = 0
$ws6.Range("G5").Value = 0
$ws6.Range("H5").Value = 0
$ws6.Range("I5").Value = 1
$ws6.Range("J5").Value = 0
$ws6.Range("K5").Value = 6
$ws6.Range("N5").Value = 0
$ws6.Range("O5").Value = 0
$ws6.Range("P5").Value = 0
$ws6.Range("Q5").Value = 5
$ws6.Range("R5").Value = 6
$ws6.Range("A6").Value = 0
$ws6.Range("B6").Value = 0
$ws6.Range("C6").Value = 0
$ws6.Range("D6").Value = 0
$ws6.Range("E6").Value = 0
$ws6.Range("F6").Value = 0
$ws6.Range("G6").Value = 0
$ws6.Range("H6").Value = 0
$ws6.Range("I6").Value = 0
$ws6.Range("J6").Value = 1
$ws6.Range("N6").Value = 0
$ws6.Range("O6").Value = 0
$ws6.Range("Q6").Value = 7
$ws6.Range("R6").Value = 18
$ws6.Range("A7").Value = 0
$ws6.Range("B7").Value = 0
$ws6.Range("C7").Value = 0
$ws6.Range("D7").Value = 0
$ws6.Range("E7").Value = 0
$ws6.Range("G7").Value = 0
$ws6.Range("N7").Value = 0
$ws6.Range("O7").Value = 0
$ws6.Range("A8").Value = 0
$ws6.Range("B8").Value = 0
$ws6.Range("C8").Value = 0
$ws6.Range("D8").Value = 0
$ws6.Range("E8").Value = 0
$ws6.Range("N8").Value = 1
$ws6.Range("A9").Value = 0
$ws6.Range("B9").Value = 0
$ws6.Range("C9").Value = 0
$ws6.Range("D9").Value = 0
$ws6.Range("E9").Value = 0
$ws6.Range("N9").Value = 1
$ws6.Range("A10").Value = 0
$ws6.Range("B10").Value = 0
$ws6.Range("C10").Value = 0
$ws6.Range("D10").Value = 0
$ws6.Range("E10").Value = 0
$ws6.Range("N10").Value = 0
$ws6.Range("A11").Value = 0
$ws6.Range("B11").Value = 0
$ws6.Range("C11").Value = 0
$ws6.Range("D11").Value = 0
$ws6.Range("E11").Value = 0
$ws6.Range("N11").Value = 0
$ws6.Range("A12").Value = 0
$ws6.Range("C12").Value = 0
$ws6.Range("D12").Value = 0
# 3) Populate new data rows (2-12) on "Unexecuted Units"
$ws7 = $wb.Worksheets("Unexecuted Units")
$ws7.Range("A2").Value = 0
$ws7.Range("B2").Value = 0
$ws7.Range("C2").Value = 0
$ws7.Range("D2").Value = 0
$ws7.Range("E2").Value = 0
$ws7.Range("F2").Value = 0
$ws7.Range("G2").Value = 0
$ws7.Range("H2").Value = 0
$ws7.Range("I2").Value = 0
$ws7.Range("J2").Value = 0
$ws7.Range("K2").Value = 0
$ws7.Range("L2").Value = 1022
$ws7.Range("N2").Value = 0
$ws7.Range("O2").Value = 0
$ws7.Range("P2").Value = 0
$ws7.Range("Q2").Value = 9
$ws7.Range("R2").Value = 16
$ws7.Range("A3").Value = 0
$ws7.Range("B3").Value = 0
$ws7.Range("C3").Value = 0
$ws7.Range("D3").Value = 0
$ws7.Range("E3").Value = 0
$ws7.Range("F3").Value = 0
$ws7.Range("G3").Value = 0
$ws7.Range("H3").Value = 0
$ws7.Range("I3").Value = 0
$ws7.Range("J3").Value = 0
$ws7.Range("K3").Value = 0
$ws7.Range("N3").Value = 0
$ws7.Range("O3").Value = 0
$ws7.Range("P3").Value = 0
$ws7.Range("Q3").Value = 0
$ws7.Range("R3").Value = 7
$ws7.Range("A4").Value = 0
$ws7.Range("B4").Value = 0
$ws7.Range("C4").Value = 0
$ws7.Range("D4").Value = 0
$ws7.Range("E4").Value = 0
$ws7.Range("F4").Value = 0
$ws7.Range("G4").Value = 0
$ws7.Range("H4").Value = 0
$ws7.Range("I4").Value = 0
$ws7.Range("J4").Value = 0
$ws7.Range("K4").Value = 263
$ws7.Range("N4").Value = 0
$ws7.Range("O4").Value = 0
$ws7.Range("P4").Value = 0
$ws7.Range("Q4").Value = 0
$ws7.Range("R4").Value = 0
$ws7.Range("A5").Value = 0
$ws7.Range("B5").Value = 0
$ws7.Range("C5").Value = 0
$ws7.Range("D5").Value = 0
$ws7.Range("E5").Value = 0
$ws7.Range("F5").Value = 0
$ws7.Range("G5").Value = 0
$ws7.Range("H5").Value = 0
$ws7.Range("I5").Value = 0
$ws7.Range("J5").Value = 0
$ws7.Range("K5").Value = 2
$ws7.Range("N5").Value = 0
$ws7.Range("O5").Value = 0
$ws7.Range("P5").Value = 0
$ws7.Range("Q5").Value = 5
$ws7.Range("R5").Value = 6
$ws7.Range("A6").Value = 0
$ws7.Range("B6").Value = 0
$ws7.Range("C6").Value = 0
$ws7.Range("D6").Value = 0
$ws7.Range("E6").Value = 0
$ws7.Range("F6").Value = 0
$ws7.Range("G6").Value = 0
$ws7.Range("H6").Value = 0
$ws7.Range("I6").Value = 0
$ws7.Range("J6").Value = 0
$ws7.Range("N6").Value = 0
$ws7.Range("O6").Value = 0
$ws7.Range("Q6").Value = 7
$ws7.Range("R6").Value = 18
$ws7.Range("A7").Value = 0
$ws7.Range("B7").Value = 0
$ws7.Range("C7").Value = 0
$ws7.Range("D7").Value = 0
$ws7.Range("E7").Value = 0
$ws7.Range("G7").Value = 0
$ws7.Range("N7").Value = 0
$ws7.Range("O7").Value = 0
$ws7.Range("A8").Value = 0
$ws7.Range("B8").Value = 0
$ws7.Range("C8").Value = 0
$ws7.Range("D8").Value = 0
$ws7.Range("E8").Value = 0
$ws7.Range("N8").Value = 1
$ws7.Range("A9").Value = 0
$ws7.Range("B9").Value = 0
$ws7.Range("C9").Value = 0
$ws7.Range("D9").Value = 0
$ws7.Range("E9").Value = 0
$ws7.Range("N9").Value = 1
$ws7.Range("A10").Value = 0
$ws7.Range("B10").Value = 0
$ws7.Range("C10").Value = 0
$ws7.Range("D10").Value = 0
$ws7.Range("E10").Value = 0
$ws7.Range("N10").Value = 0
$ws7.Range("A11").Value = 0
$ws7.Range("B11").Value = 0
$ws7.Range("C11").Value = 0
$ws7.Range("D11").Value = 0
$ws7.Range("E11").Value = 0
$ws7.Range("N11").Value = 0
$ws7.Range("A12").Value = 0
$ws7.Range("C12").Value = 0
$ws7.Range("D12").Value = 0
# 4) Fix up selection / active sheet state.
# "Unexecuted Units" keeps a non-active selection...
$ws7.Range("R2").Select()
# ...while "Incomplete Units" ends up the active/selected sheet (this also
# clears tabSelected on whichever sheet - "Success Rate" - was active before).
$ws6.Range("P18").Select()
